$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Insert a brand-new column before column A. This shifts every existing
#    column (values, styles, formulas, data validations) one slot to the
#    right -- exactly what the diff shows (A -> B, B -> C, ... AB -> AC).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).Insert()

# ---------------------------------------------------------------------------
# 2. The newly inserted column has no formatting yet. Clone the formatting
#    of the (now shifted) former column A -- currently column B -- onto the
#    new column A so header/data/styling look consistent with the rest of
#    the table.
# ---------------------------------------------------------------------------
$ws.Range("B1:B12").Copy()
$ws.Range("A1:A12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Match the new column's width (23 chars). The engine stores width in the
#    OOXML `width` attribute as ColumnWidth + 5/6, so back that offset out.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23 - 5/6

# ---------------------------------------------------------------------------
# 4. Populate the new "INDEX (DO NOT MODIFY)" column.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"

$indexValues = @(43, 44, 45, 46, 47, 48, 49, 146, 174, 193, 245)
for ($i = 0; $i -lt $indexValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $indexValues[$i]
}

# ---------------------------------------------------------------------------
# 5. Upper-case the rest of the header row (now columns B..AB -- the former
#    A1..AA1). The trailing status column (now AC1, formerly AB1) is left
#    untouched ("Status as of July 11, 2025" keeps its original casing).
# ---------------------------------------------------------------------------
for ($c = 2; $c -le 28; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = $cell.Value().ToUpper()
}
